# Upload new version with timestamp
# Inserts a new "OMEGA P SYRUP 120 ML" shortage-list row (between the existing
# "OFRAMAX 1 GM I.M. VIAL" row and "OTAL EAR DROPS 5 ML" row), renumbers the
# serial column for the rows that shift down, bumps the running total, and
# refreshes the generated timestamp in the report footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds "OTAL EAR DROPS 5 ML" (and everything below it, down
# to the totals/footer rows, shifts down by one row). Insert a fresh row at
# 23 - Excel shifts rows 23..34 down to 24..35, carrying over formatting.
$ws.Rows.Item(23).Insert(-4121)

# Fill in the new row with the new item's data. The "حد الطلب" / price /
# sale-price columns are stored as text in this report (e.g. "71.0000", not
# 71), so force them through the "@" text format while assigning - otherwise
# Excel auto-converts the numeric-looking strings to real numbers.
$ws.Range("A23").Value = 17
$ws.Range("C23").Value = "OMEGA P SYRUP 120 ML"
$ws.Range("H23").Value = "0:0"

$Lfmt = $ws.Range("L23").NumberFormat
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = "1"
$ws.Range("L23").NumberFormat = $Lfmt

$Nfmt = $ws.Range("N23").NumberFormat
$ws.Range("N23").NumberFormat = "@"
$ws.Range("N23").Value = "125.00"
$ws.Range("N23").NumberFormat = $Nfmt

$Pfmt = $ws.Range("P23").NumberFormat
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = "125.0000"
$ws.Range("P23").NumberFormat = $Pfmt

$ws.Range("Q23").Value = "1:0"

# Renumber the serial ("م") column for the rows that moved down one slot.
for ($r = 24; $r -le 32; $r++) {
    $ws.Range("A$r").Value = $r - 6
}

# The grand total (was 1221.26) grows by the new item's sale price (125.00).
$ws.Range("P33").Value = 1346.26

# Refresh the generated-on timestamp shown in the report footer.
$ws.Range("A34").Value = "Saturday, 2 August, 2025 12:34 PM"
